$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for specific rows per the diff
$ws.Range("F2").Value = 8
$ws.Range("F4").Value = 7
$ws.Range("F5").Value = -9
$ws.Range("F6").Value = 1
$ws.Range("F10").Value = -1
